$wb = $excel.ActiveWorkbook

# --- Sheet2 ("SoFCtMbCtPR"): update hard coal share value, adjust selection ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = 2.5
$null = $ws2.Range("B3").Select()

# --- Sheet1 ("About"): append the Notes section ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A5").Value = "Notes:"
$ws1.Range("A5").Font.Bold = $true

$ws1.Range("B5").Value = "In the US, many coal plants are subject to rules requiring them "
$ws1.Range("B6").Value = "to retrofit to meet enviromental guidelines. This requires"
$ws1.Range("B7").Value = "a one time investment decision for plant owners. Because we don't"
$ws1.Range("B8").Value = "track individual plants in the model, we calibrate the share of forward"
$ws1.Range("B9").Value = "costs that must be recovered to represent the additional revenue that is needed to "
$ws1.Range("B10").Value = "save and pay for these one time investments and apply this across the distribution"
$ws1.Range("B11").Value = "of plant types. Calibration is done by comparing model results against other sources,"
$ws1.Range("B12").Value = "including Rhodium's ClimateDeck and EIA's Annual Energy Outlook and Electric "
$ws1.Range("B13").Value = "Power Monthly."

# Re-activate the "About" sheet last so it keeps the tab selection, and
# restore its selected cell.
$null = $ws1.Activate()
$null = $ws1.Range("D31").Select()
